# Update the cryptocurrency price/volume table with refreshed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a text value to a cell. If the text would otherwise be
# auto-interpreted by Excel as a number (losing formatting such as trailing
# zeroes, e.g. '1.000' becoming 1), prefix it with an apostrophe so Excel
# keeps it as literal text, exactly as it was stored in the source data.
function Set-TextValue($cellRef, $text) {
    if ($text -match '^\s*[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?\s*$') {
        $ws.Range($cellRef).Value = "'" + $text
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

Set-TextValue 'D2' '29.442.33'
Set-TextValue 'E2' '  +0.40%  '
Set-TextValue 'D3' '1.849.71'
Set-TextValue 'E3' '  +0.42%  '
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  +0.14%  '
Set-TextValue 'D5' '240.69'
Set-TextValue 'E5' '  +0.73%  '
Set-TextValue 'D6' '0.6295'
Set-TextValue 'E6' '  -0.04%  '
Set-TextValue 'D7' '1.000'
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'D8' '0.07704'
Set-TextValue 'E9' '  -0.53%  '
Set-TextValue 'D10' '24.69'
Set-TextValue 'E10' '  +0.83%  '
Set-TextValue 'D11' '0.07740'
Set-TextValue 'E11' '  +0.68%  '
Set-TextValue 'D12' '1.881.17'
Set-TextValue 'E12' '  +1.60%  '
Set-TextValue 'E13' '  +1.24%  '
Set-TextValue 'E14' '  +3.95%  '
Set-TextValue 'D15' '0.6790'
Set-TextValue 'D16' '83.65'
Set-TextValue 'E16' '  +0.70%  '
Set-TextValue 'D17' '2.145.97'
Set-TextValue 'E17' '  +1.43%  '
Set-TextValue 'D18' '6.196'
Set-TextValue 'E18' '  +0.84%  '
Set-TextValue 'D19' '29.471.74'
Set-TextValue 'E19' '  +0.34%  '
Set-TextValue 'D20' '228.33'
Set-TextValue 'E20' '  -0.13%  '
Set-TextValue 'E21' '  +0.27%  '
Set-TextValue 'E22' '  +0.04%  '
Set-TextValue 'D23' '7.439'
Set-TextValue 'E23' '  -0.12%  '
Set-TextValue 'D24' '1.000'
Set-TextValue 'E24' '  +0.02%  '
Set-TextValue 'D25' '157.59'
Set-TextValue 'E25' '  +0.83%  '
Set-TextValue 'D26' '0.1379'
Set-TextValue 'E26' '  -1.04%  '
Set-TextValue 'D27' '8.415'
Set-TextValue 'E27' '  +0.73%  '
Set-TextValue 'D28' '17.69'
Set-TextValue 'E28' '  +0.47%  '
Set-TextValue 'D29' '1.344'
Set-TextValue 'E29' '  +5.80%  '
Set-TextValue 'D30' '1.467'
Set-TextValue 'E30' '  +0.54%  '
Set-TextValue 'D31' '0.05679'
Set-TextValue 'E32' '  +0.44%  '
Set-TextValue 'E33' '  +0.45%  '
Set-TextValue 'D34' '1.846'
Set-TextValue 'D35' '1.163'
Set-TextValue 'E35' '  +0.82%  '
Set-TextValue 'D36' '0.7036'
Set-TextValue 'E36' '  -0.75%  '
Set-TextValue 'E37' '  -0.19%  '
Set-TextValue 'D38' '2.780'
Set-TextValue 'E38' '  +0.46%  '
Set-TextValue 'D39' '0.01791'
Set-TextValue 'E39' '  -0.94%  '
Set-TextValue 'D40' '1.220.33'
Set-TextValue 'E40' '  -1.65%  '
Set-TextValue 'D41' '6.552'
Set-TextValue 'E41' '  +5.08%  '
Set-TextValue 'D42' '0.9047'
Set-TextValue 'E42' '  +0.42%  '
Set-TextValue 'E43' '  +0.15%  '
Set-TextValue 'D44' '101.86'
Set-TextValue 'E44' '  +0.07%  '
Set-TextValue 'D45' '66.30'
Set-TextValue 'E45' '  +1.31%  '
Set-TextValue 'B46' 'Aptos'
Set-TextValue 'C46' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D46' '7.141'
Set-TextValue 'E46' '  +0.56%  '
Set-TextValue 'B47' 'BabyDogeCoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D47' '0.00000000119'
Set-TextValue 'E47' '  +0.49%  '
Set-TextValue 'D48' '0.4019'
Set-TextValue 'D49' '9.011'
Set-TextValue 'E49' '  +0.85%  '
Set-TextValue 'E50' '  +0.55%  '
Set-TextValue 'D51' '0.1146'
Set-TextValue 'E51' '  +2.19%  '
